$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 3 through 47 (A, B, C values). Row 3-12 already existed (only C changes),
# rows 13-47 are newly added.
$data = @(
    @(1, 1, 0),
    @(2, 2, 1),
    @(3, 3, 1),
    @(4, 4, 2),
    @(5, 5, 3),
    @(6, 6, 3),
    @(7, 7, 4),
    @(8, 8, 5),
    @(9, 9, 5),
    @(10, 10, 6),
    @(11, 11, 7),
    @(12, 12, 7),
    @(13, 13, 8),
    @(14, 14, 9),
    @(15, 15, 9),
    @(16, 16, 10),
    @(17, 17, 10),
    @(18, 18, 11),
    @(19, 19, 12),
    @(20, 20, 12),
    @(21, 21, 13),
    @(22, 22, 14),
    @(23, 23, 14),
    @(24, 24, 15),
    @(25, 25, 16),
    @(26, 26, 16),
    @(27, 27, 17),
    @(28, 28, 18),
    @(29, 29, 18),
    @(30, 30, 19),
    @(31, 31, 19),
    @(32, 32, 20),
    @(33, 33, 21),
    @(34, 34, 21),
    @(35, 35, 22),
    @(36, 36, 23),
    @(37, 37, 23),
    @(38, 38, 24),
    @(39, 39, 25),
    @(40, 40, 25),
    @(41, 41, 26),
    @(42, 42, 27),
    @(43, 43, 27),
    @(44, 44, 28),
    @(45, 45, 29)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
